# Applies the cryptos.xlsx price/volume refresh + three coin re-ranks described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row re-orders (Coin name + Link swap between two existing rows) ---
$rowSwaps = @(
    @{ Row = 23; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Row = 24; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Row = 26; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Row = 27; B = 'LEO'; C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Row = 44; B = 'Fetch.AI'; C = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' },
    @{ Row = 45; B = 'ThetaToken'; C = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta' },
    @{ Row = 48; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Row = 49; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Row = 50; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Row = 51; B = 'ApeXProtocol'; C = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
)

foreach ($item in $rowSwaps) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
}

# --- Price (D) and Volume(1h) (E) refresh for every data row (2-51) ---
# NumberFormat is forced to Text ("@") before assignment so that numeric-looking
# price strings (e.g. "1.00", "9.00") are preserved verbatim instead of Excel
# normalizing them into real numbers (which would drop the literal formatting).
$priceVolume = @(
    @{ Row = 2; D = '63.009.56'; E = '  -7.93%  ' },
    @{ Row = 3; D = '3.243.37'; E = '  -10.14%  ' },
    @{ Row = 4; D = '0.998'; E = '  -0.52%  ' },
    @{ Row = 5; D = '174.89'; E = '  -15.27%  ' },
    @{ Row = 6; D = '508.39'; E = '  -10.96%  ' },
    @{ Row = 7; D = '0.587'; E = '  -4.45%  ' },
    @{ Row = 8; E = '  +0.09%  ' },
    @{ Row = 9; D = '3.241.73'; E = '  -10.06%  ' },
    @{ Row = 10; D = '0.608'; E = '  -11.42%  ' },
    @{ Row = 11; D = '56.23'; E = '  -11.86%  ' },
    @{ Row = 12; D = '0.128'; E = '  -13.73%  ' },
    @{ Row = 13; D = '0.0000250'; E = '  -11.54%  ' },
    @{ Row = 14; D = '9.00'; E = '  -12.97%  ' },
    @{ Row = 15; D = '3.731.65'; E = '  -10.82%  ' },
    @{ Row = 16; D = '0.118'; E = '  -6.75%  ' },
    @{ Row = 17; D = '3.224.41'; E = '  -10.35%  ' },
    @{ Row = 18; D = '62.780.31'; E = '  -8.04%  ' },
    @{ Row = 19; D = '17.10'; E = '  -11.66%  ' },
    @{ Row = 20; D = '10.76'; E = '  -12.49%  ' },
    @{ Row = 21; D = '0.930'; E = '  -13.26%  ' },
    @{ Row = 22; D = '365.06'; E = '  -10.24%  ' },
    @{ Row = 23; D = '10.92'; E = '  -12.30%  ' },
    @{ Row = 24; D = '78.67'; E = '  -7.50%  ' },
    @{ Row = 25; D = '3.59'; E = '  -14.36%  ' },
    @{ Row = 26; D = '3.75'; E = '  -2.87%  ' },
    @{ Row = 27; D = '5.92'; E = '  -3.53%  ' },
    @{ Row = 28; D = '2.59'; E = '  -11.03%  ' },
    @{ Row = 29; D = '11.12'; E = '  -11.48%  ' },
    @{ Row = 30; D = '8.14'; E = '  -12.02%  ' },
    @{ Row = 31; D = '641.26'; E = '  -9.82%  ' },
    @{ Row = 32; D = '27.91'; E = '  -12.07%  ' },
    @{ Row = 33; D = '6.53'; E = '  -15.25%  ' },
    @{ Row = 34; D = '11.00'; E = '  -9.91%  ' },
    @{ Row = 35; D = '57.52'; E = '  -9.71%  ' },
    @{ Row = 36; D = '0.102'; E = '  -10.43%  ' },
    @{ Row = 37; E = '  +0.08%  ' },
    @{ Row = 38; D = '35.33'; E = '  -16.14%  ' },
    @{ Row = 39; D = '0.374'; E = '  -10.55%  ' },
    @{ Row = 40; D = '0.994'; E = '  -0.47%  ' },
    @{ Row = 41; D = '0.122'; E = '  -8.19%  ' },
    @{ Row = 42; D = '2.848.29'; E = '  -10.56%  ' },
    @{ Row = 43; D = '0.0₃0643'; E = '  -15.15%  ' },
    @{ Row = 44; D = '2.38'; E = '  -10.70%  ' },
    @{ Row = 45; D = '2.62'; E = '  -20.53%  ' },
    @{ Row = 46; D = '2.55'; E = '  -8.73%  ' },
    @{ Row = 47; D = '2.77'; E = '  +2.34%  ' },
    @{ Row = 48; D = '25.34'; E = '  +13.73%  ' },
    @{ Row = 49; D = '0.0375'; E = '  -9.98%  ' },
    @{ Row = 50; D = '0.122'; E = '  -7.25%  ' },
    @{ Row = 51; D = '2.88'; E = '  -6.68%  ' }
)

foreach ($item in $priceVolume) {
    $r = $item.Row
    if ($item.ContainsKey("D")) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $item.E
    }
}

